$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts old row 13 "HX115" down to row 14)
$ws.Rows.Item(13).Insert()

# --- Update existing rows 2-11 (values in columns C and D changed) ---
$ws.Cells.Item(2, 3).Value = 829300028.4920082
$ws.Cells.Item(2, 4).Value = 4.57412869797354

$ws.Cells.Item(3, 3).Value = 444337009.392413
$ws.Cells.Item(3, 4).Value = 2.450807423616483

$ws.Cells.Item(4, 4).Value = 0.107689515177348

$ws.Cells.Item(5, 3).Value = 10305322.15321375
$ws.Cells.Item(5, 4).Value = 0.05684055008245033

$ws.Cells.Item(6, 4).Value = 0.01282226808070471

$ws.Cells.Item(7, 4).Value = 0.0000329563936755097

$ws.Cells.Item(8, 3).Value = -2019186.584710486
$ws.Cells.Item(8, 4).Value = -0.01113712647578478

$ws.Cells.Item(9, 4).Value = -0.05295859353664479

$ws.Cells.Item(10, 3).Value = -18190636.23184341
$ws.Cells.Item(10, 4).Value = -0.1003331826405141

$ws.Cells.Item(11, 4).Value = -0.2942570636440101

# --- Row 12: was FLASH107/Flash, now becomes HXN/HeatExchangerNetwork ---
$ws.Cells.Item(12, 1).Value = "HXN"
$ws.Cells.Item(12, 2).Value = "HeatExchangerNetwork"
$ws.Cells.Item(12, 3).Value = -55986681.98138828
$ws.Cells.Item(12, 4).Value = -0.3088029421885569

# --- New row 13: FLASH107/Flash with new values ---
$ws.Cells.Item(13, 1).Value = "FLASH107"
$ws.Cells.Item(13, 2).Value = "Flash"
$ws.Cells.Item(13, 3).Value = -323156771.1969261
$ws.Cells.Item(13, 4).Value = -1.782419643424109

# --- Row 14 (formerly row 13, HX115): only D14 value changed ---
$ws.Cells.Item(14, 4).Value = -3.65241285941458
